$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.370.67"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "3.515.80"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'592.09"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").Value = "'134.76"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.488"
$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").Value = "'7.62"
$ws.Range("E9").Value = "  +6.92%  "

$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("E11").Value = "  +4.13%  "

$ws.Range("D12").Value = "4.115.51"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").Value = "3.517.37"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "'25.87"
$ws.Range("E16").Value = "  -2.03%  "

$ws.Range("D17").Value = "64.371.91"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").Value = "'10.01"
$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").Value = "'5.78"
$ws.Range("E19").Value = "  +3.48%  "

$ws.Range("D20").Value = "'13.58"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").Value = "'394.37"
$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("D23").Value = "3.656.80"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("D24").Value = "'74.63"
$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  +2.80%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("E29").Value = "  -1.92%  "

$ws.Range("E30").Value = "  +1.48%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("E32").Value = "  -6.57%  "

$ws.Range("E33").Value = "  +7.54%  "

$ws.Range("D34").Value = "3.546.88"
$ws.Range("E34").Value = "  +0.70%  "

$ws.Range("D36").Value = "'23.44"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").Value = "'6.98"
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("D40").Value = "'166.94"
$ws.Range("E40").Value = "  +1.48%  "

$ws.Range("E41").Value = "  +0.98%  "

$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("D43").Value = "'25.59"
$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'4.45"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").Value = "'1.19"
$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").Value = "2.399.19"
$ws.Range("E49").Value = "  -3.16%  "

$ws.Range("D50").Value = "'0.901"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("E51").Value = "  +0.28%  "
